# Weekly crime-stat refresh: new week label/date range + updated row 15-33 figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/issue number and reporting week date range ---
$ws.Range("A8").Value = "Volume 32   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  6/30/2025  Through  7/6/2025"

# --- Cells whose style/type changes (placeholder "0"/"***.*"  <->  real numbers) ---
# PasteSpecial(-4122) = xlPasteFormats copies the number-format/style only (no value).
# PasteSpecial(-4163) = xlPasteAll, run a second time from the same placeholder cell,
# is used when the destination also needs the exact shared placeholder text/value.
$ws.Range("F14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F14").Copy()
$ws.Range("F15").PasteSpecial(-4163)
$ws.Range("G14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("G15").PasteSpecial(-4163)
$ws.Range("H14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("H15").PasteSpecial(-4163)
$ws.Range("C16").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 1
$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("F14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F14").Copy()
$ws.Range("F27").PasteSpecial(-4163)
$ws.Range("G14").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("G27").PasteSpecial(-4163)
$ws.Range("H14").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("H27").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("F16").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F31").Value = 1
$ws.Range("D17").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = 4
$ws.Range("E17").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = -100
$ws.Range("G16").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("G33").Value = 4
$ws.Range("H16").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").Value = -100

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("C16").Value = 4
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 166.666666666667
$ws.Range("I16").Value = 61
$ws.Range("K16").Value = 29.787234042553
$ws.Range("L16").Value = -10.294117647058
$ws.Range("M16").Value = -17.567567567567
$ws.Range("N16").Value = -84.478371501272
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -6.25
$ws.Range("I17").Value = 112
$ws.Range("J17").Value = 123
$ws.Range("K17").Value = -8.943089430894
$ws.Range("L17").Value = 4.672897196261
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = -0.884955752212
$ws.Range("E18").Value = -66.666666666666
$ws.Range("I18").Value = 47
$ws.Range("J18").Value = 75
$ws.Range("K18").Value = -37.333333333333
$ws.Range("L18").Value = -41.25
$ws.Range("M18").Value = 27.027027027027
$ws.Range("N18").Value = -74.175824175824
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -36.170212765957
$ws.Range("I19").Value = 207
$ws.Range("J19").Value = 222
$ws.Range("K19").Value = -6.756756756756
$ws.Range("L19").Value = -31.229235880398
$ws.Range("M19").Value = 56.818181818181
$ws.Range("N19").Value = -9.606986899563
$ws.Range("C20").Value = 4
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 17
$ws.Range("K20").Value = -43.333333333333
$ws.Range("L20").Value = -41.379310344827
$ws.Range("M20").Value = -37.037037037037
$ws.Range("N20").Value = -91.145833333333
$ws.Range("C21").Value = 22
$ws.Range("E21").Value = 37.5
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = -24.691358024691
$ws.Range("I21").Value = 452
$ws.Range("J21").Value = 503
$ws.Range("K21").Value = -10.139165009940
$ws.Range("L21").Value = -23.648648648648
$ws.Range("M21").Value = 27.323943661971
$ws.Range("N21").Value = -59.714795008912
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("M22").Value = 157.142857142857
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 9.090909090909
$ws.Range("I23").Value = 86
$ws.Range("J23").Value = 111
$ws.Range("K23").Value = -22.522522522522
$ws.Range("L23").Value = 10.256410256410
$ws.Range("M23").Value = 11.688311688311
$ws.Range("C24").Value = 28
$ws.Range("E24").Value = -6.666666666666
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 122
$ws.Range("H24").Value = 3.278688524590
$ws.Range("I24").Value = 599
$ws.Range("J24").Value = 675
$ws.Range("K24").Value = -11.259259259259
$ws.Range("L24").Value = -3.231017770597
$ws.Range("M24").Value = 56.396866840731
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -26.923076923076
$ws.Range("F25").Value = 73
$ws.Range("G25").Value = 93
$ws.Range("H25").Value = -21.505376344086
$ws.Range("I25").Value = 339
$ws.Range("J25").Value = 459
$ws.Range("K25").Value = -26.143790849673
$ws.Range("L25").Value = -0.586510263929
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 50
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 34.482758620689
$ws.Range("I26").Value = 208
$ws.Range("J26").Value = 222
$ws.Range("K26").Value = -6.306306306306
$ws.Range("L26").Value = -10.729613733905
$ws.Range("M26").Value = 17.514124293785
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = 0
$ws.Range("L29").Value = -66.666666666666
$ws.Range("M29").Value = -50
$ws.Range("L30").Value = -60
$ws.Range("M30").Value = -33.333333333333
$ws.Range("I31").Value = 4
$ws.Range("K31").Value = -55.555555555555
$ws.Range("L31").Value = 33.333333333333
$ws.Range("J33").Value = 5
